# Auto-generated-assisted edit script for cfb_weather.xlsx
# Updates odds/weather derived metrics (columns S,T,V,Y,Z,AB) refreshed by the
# scraper re-run, removes stale U29/V29/Z29 cells, and refreshes the shared
# "Timestamp" column (AE) string for all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FBS")

# --- Update changed numeric values ---
$ws.Range("V4").Value = -27.5
$ws.Range("Z4").Value = 0
$ws.Range("V6").Value = -6
$ws.Range("Z6").Value = 0.5
$ws.Range("S9").Value = 49.5
$ws.Range("V9").Value = 9
$ws.Range("Y9").Value = -0.07476635514018691
$ws.Range("AB9").Value = -0.1236037414153002
$ws.Range("Z9").Value = -2.5
$ws.Range("V13").Value = 11
$ws.Range("Z13").Value = -1
$ws.Range("V16").Value = -7
$ws.Range("Z16").Value = -0.5
$ws.Range("V17").Value = 3
$ws.Range("Z17").Value = -1
$ws.Range("V18").Value = 9
$ws.Range("Z18").Value = -1.5
$ws.Range("T21").Value = -110
$ws.Range("V23").Value = -7
$ws.Range("Z23").Value = 0
$ws.Range("V24").Value = -3
$ws.Range("Z24").Value = -0.5
$ws.Range("S26").Value = 34.5
$ws.Range("V26").Value = 9.5
$ws.Range("T26").Value = -110
$ws.Range("Z26").Value = -2
$ws.Range("AB26").Value = -0.3365384615384616
$ws.Range("Y26").Value = -0.0547945205479452
$ws.Range("V28").Value = 9.5
$ws.Range("Z28").Value = -2
$ws.Range("V30").Value = -2.5
$ws.Range("Z30").Value = 0.5
$ws.Range("V32").Value = 13.5
$ws.Range("Z32").Value = -2
$ws.Range("S33").Value = 49.5
$ws.Range("T33").Value = -110
$ws.Range("V33").Value = -16
$ws.Range("AB33").Value = -0.1277533039647577
$ws.Range("Z33").Value = 2.5
$ws.Range("Y33").Value = -0.0198019801980198
$ws.Range("V34").Value = -19.5
$ws.Range("Z34").Value = -0.5
$ws.Range("V36").Value = 14
$ws.Range("Z36").Value = 0.5
$ws.Range("S37").Value = 54.5
$ws.Range("Y37").Value = -0.01801801801801802
$ws.Range("AB37").Value = -0.03964757709251102
$ws.Range("V39").Value = 14
$ws.Range("Z39").Value = 0
$ws.Range("V40").Value = 3
$ws.Range("Z40").Value = 0.5
$ws.Range("T45").Value = -110
$ws.Range("V46").Value = -7.5
$ws.Range("Z46").Value = -1
$ws.Range("S47").Value = 55.5
$ws.Range("T47").Value = -118
$ws.Range("Y47").Value = 0
$ws.Range("AB47").Value = -0.03896103896103896

# --- Clear cells that no longer have data in the refreshed pull ---
$ws.Range("U29").ClearContents()
$ws.Range("V29").ClearContents()
$ws.Range("Z29").ClearContents()

# --- Refresh the pull Timestamp for every data row (shared by all rows) ---
$ws.Range("AE2:AE47").Value = "2024-09-30T16:21:36.611528"

